$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("Schedule") ----
$ws1 = $wb.Worksheets.Item("Schedule")

# Update existing row 2 values
$ws1.Cells.Item(2, 1).Value = 46041.16666666666
$ws1.Cells.Item(2, 2).Value = 46041.66666666666
$ws1.Cells.Item(2, 3).Value = 12
$ws1.Cells.Item(2, 4).Value = 45.36
$ws1.Cells.Item(2, 5).Value = 341.6364315
$ws1.Cells.Item(2, 6).Value = 7.53166736111111

# Add new row 3, copying formatting (date styles) from row 2 first
$ws1.Range("A2:F2").Copy() | Out-Null
$ws1.Range("A3:F3").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(3, 1).Value = 46041.83333333334
$ws1.Cells.Item(3, 2).Value = 46042
$ws1.Cells.Item(3, 3).Value = 4
$ws1.Cells.Item(3, 4).Value = 15.12
$ws1.Cells.Item(3, 5).Value = 484.6414267499999
$ws1.Cells.Item(3, 6).Value = 32.05300441468254

# ---- Sheet2 ("Detailed") ----
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Cells.Item(2, 1).Value = 46041.02083333334
$ws2.Cells.Item(2, 2).Value = 57.06003
$ws2.Cells.Item(2, 3).Value = "historical"
$ws2.Cells.Item(2, 5).Value = "OFF"
$ws2.Cells.Item(3, 1).Value = 46041.04166666666
$ws2.Cells.Item(3, 2).Value = 57.06003
$ws2.Cells.Item(3, 3).Value = "historical"
$ws2.Cells.Item(3, 5).Value = "OFF"
$ws2.Cells.Item(4, 1).Value = 46041.0625
$ws2.Cells.Item(4, 2).Value = 56.98
$ws2.Cells.Item(4, 3).Value = "forecast"
$ws2.Cells.Item(4, 5).Value = "OFF"
$ws2.Cells.Item(5, 1).Value = 46041.08333333334
$ws2.Cells.Item(5, 2).Value = 47.65133
$ws2.Cells.Item(5, 3).Value = "forecast"
$ws2.Cells.Item(5, 5).Value = "OFF"
$ws2.Cells.Item(6, 1).Value = 46041.10416666666
$ws2.Cells.Item(6, 2).Value = 47.99544
$ws2.Cells.Item(6, 3).Value = "forecast"
$ws2.Cells.Item(6, 5).Value = "OFF"
$ws2.Cells.Item(7, 1).Value = 46041.125
$ws2.Cells.Item(7, 2).Value = 48.31738
$ws2.Cells.Item(7, 3).Value = "forecast"
$ws2.Cells.Item(7, 5).Value = "OFF"
$ws2.Cells.Item(8, 1).Value = 46041.14583333334
$ws2.Cells.Item(8, 2).Value = 57.06003
$ws2.Cells.Item(8, 3).Value = "forecast"
$ws2.Cells.Item(8, 5).Value = "OFF"
$ws2.Cells.Item(9, 1).Value = 46041.16666666666
$ws2.Cells.Item(9, 2).Value = 57.06003
$ws2.Cells.Item(9, 3).Value = "forecast"
$ws2.Cells.Item(9, 5).Value = "ON"
$ws2.Cells.Item(10, 1).Value = 46041.1875
$ws2.Cells.Item(10, 2).Value = 58.82614
$ws2.Cells.Item(10, 3).Value = "forecast"
$ws2.Cells.Item(10, 5).Value = "ON"
$ws2.Cells.Item(11, 1).Value = 46041.20833333334
$ws2.Cells.Item(11, 2).Value = 57.83627
$ws2.Cells.Item(11, 3).Value = "forecast"
$ws2.Cells.Item(11, 5).Value = "ON"
$ws2.Cells.Item(12, 1).Value = 46041.22916666666
$ws2.Cells.Item(12, 2).Value = 60.45412
$ws2.Cells.Item(12, 3).Value = "forecast"
$ws2.Cells.Item(12, 5).Value = "ON"
$ws2.Cells.Item(13, 1).Value = 46041.25
$ws2.Cells.Item(13, 2).Value = 60.20735
$ws2.Cells.Item(13, 3).Value = "forecast"
$ws2.Cells.Item(13, 5).Value = "ON"
$ws2.Cells.Item(14, 1).Value = 46041.27083333334
$ws2.Cells.Item(14, 2).Value = 57.06003
$ws2.Cells.Item(14, 3).Value = "forecast"
$ws2.Cells.Item(14, 5).Value = "ON"
$ws2.Cells.Item(15, 1).Value = 46041.29166666666
$ws2.Cells.Item(15, 2).Value = 32.50105
$ws2.Cells.Item(15, 3).Value = "forecast"
$ws2.Cells.Item(15, 5).Value = "ON"
$ws2.Cells.Item(16, 1).Value = 46041.3125
$ws2.Cells.Item(16, 2).Value = 1.302
$ws2.Cells.Item(16, 3).Value = "forecast"
$ws2.Cells.Item(16, 5).Value = "ON"
$ws2.Cells.Item(17, 1).Value = 46041.33333333334
$ws2.Cells.Item(17, 2).Value = 0.51
$ws2.Cells.Item(17, 3).Value = "forecast"
$ws2.Cells.Item(17, 5).Value = "ON"
$ws2.Cells.Item(18, 1).Value = 46041.35416666666
$ws2.Cells.Item(18, 2).Value = 6.66416
$ws2.Cells.Item(18, 3).Value = "forecast"
$ws2.Cells.Item(18, 5).Value = "ON"
$ws2.Cells.Item(19, 1).Value = 46041.375
$ws2.Cells.Item(19, 2).Value = 0.51
$ws2.Cells.Item(19, 3).Value = "forecast"
$ws2.Cells.Item(19, 5).Value = "ON"
$ws2.Cells.Item(20, 1).Value = 46041.39583333334
$ws2.Cells.Item(20, 2).Value = -0.9374400000000001
$ws2.Cells.Item(20, 3).Value = "forecast"
$ws2.Cells.Item(20, 5).Value = "ON"
$ws2.Cells.Item(21, 1).Value = 46041.41666666666
$ws2.Cells.Item(21, 2).Value = -5.50985
$ws2.Cells.Item(21, 3).Value = "forecast"
$ws2.Cells.Item(21, 5).Value = "ON"
$ws2.Cells.Item(22, 1).Value = 46041.4375
$ws2.Cells.Item(22, 2).Value = -5.92186
$ws2.Cells.Item(22, 3).Value = "forecast"
$ws2.Cells.Item(22, 5).Value = "ON"
$ws2.Cells.Item(23, 1).Value = 46041.45833333334
$ws2.Cells.Item(23, 2).Value = -0.93203
$ws2.Cells.Item(23, 3).Value = "forecast"
$ws2.Cells.Item(23, 5).Value = "ON"
$ws2.Cells.Item(24, 1).Value = 46041.47916666666
$ws2.Cells.Item(24, 2).Value = -4.93017
$ws2.Cells.Item(24, 3).Value = "forecast"
$ws2.Cells.Item(24, 5).Value = "ON"
$ws2.Cells.Item(25, 1).Value = 46041.5
$ws2.Cells.Item(25, 2).Value = -5.01
$ws2.Cells.Item(25, 3).Value = "forecast"
$ws2.Cells.Item(25, 5).Value = "ON"
$ws2.Cells.Item(26, 1).Value = 46041.52083333334
$ws2.Cells.Item(26, 2).Value = -5.75827
$ws2.Cells.Item(26, 3).Value = "forecast"
$ws2.Cells.Item(26, 5).Value = "ON"
$ws2.Cells.Item(27, 1).Value = 46041.54166666666
$ws2.Cells.Item(27, 2).Value = -5.50985
$ws2.Cells.Item(27, 3).Value = "forecast"
$ws2.Cells.Item(27, 5).Value = "ON"
$ws2.Cells.Item(28, 1).Value = 46041.5625
$ws2.Cells.Item(28, 2).Value = -5.01
$ws2.Cells.Item(28, 3).Value = "forecast"
$ws2.Cells.Item(28, 5).Value = "ON"
$ws2.Cells.Item(29, 1).Value = 46041.58333333334
$ws2.Cells.Item(29, 2).Value = -2.61261
$ws2.Cells.Item(29, 3).Value = "forecast"
$ws2.Cells.Item(29, 5).Value = "ON"
$ws2.Cells.Item(30, 1).Value = 46041.60416666666
$ws2.Cells.Item(30, 2).Value = -0.91299
$ws2.Cells.Item(30, 3).Value = "forecast"
$ws2.Cells.Item(30, 5).Value = "ON"
$ws2.Cells.Item(31, 1).Value = 46041.625
$ws2.Cells.Item(31, 2).Value = 0.00026
$ws2.Cells.Item(31, 3).Value = "forecast"
$ws2.Cells.Item(31, 5).Value = "ON"
$ws2.Cells.Item(32, 1).Value = 46041.64583333334
$ws2.Cells.Item(32, 2).Value = 0.51
$ws2.Cells.Item(32, 3).Value = "forecast"
$ws2.Cells.Item(32, 5).Value = "ON"
$ws2.Cells.Item(33, 1).Value = 46041.66666666666
$ws2.Cells.Item(33, 2).Value = 0.7
$ws2.Cells.Item(33, 3).Value = "forecast"
$ws2.Cells.Item(33, 5).Value = "OFF"
$ws2.Cells.Item(34, 1).Value = 46041.6875
$ws2.Cells.Item(34, 2).Value = -2.47201
$ws2.Cells.Item(34, 3).Value = "forecast"
$ws2.Cells.Item(34, 5).Value = "OFF"
$ws2.Cells.Item(35, 1).Value = 46041.70833333334
$ws2.Cells.Item(35, 2).Value = -2.54304
$ws2.Cells.Item(35, 3).Value = "forecast"
$ws2.Cells.Item(35, 5).Value = "OFF"
$ws2.Cells.Item(36, 1).Value = 46041.72916666666
$ws2.Cells.Item(36, 2).Value = 0.00957
$ws2.Cells.Item(36, 3).Value = "forecast"
$ws2.Cells.Item(36, 5).Value = "OFF"
$ws2.Cells.Item(37, 1).Value = 46041.75
$ws2.Cells.Item(37, 2).Value = 11.99698
$ws2.Cells.Item(37, 3).Value = "forecast"
$ws2.Cells.Item(37, 5).Value = "OFF"
$ws2.Cells.Item(38, 1).Value = 46041.77083333334
$ws2.Cells.Item(38, 2).Value = 52.39947
$ws2.Cells.Item(38, 3).Value = "forecast"
$ws2.Cells.Item(38, 5).Value = "OFF"
$ws2.Cells.Item(39, 1).Value = 46041.79166666666
$ws2.Cells.Item(39, 2).Value = 56.74129
$ws2.Cells.Item(39, 3).Value = "forecast"
$ws2.Cells.Item(39, 5).Value = "OFF"
$ws2.Cells.Item(40, 1).Value = 46041.8125
$ws2.Cells.Item(40, 2).Value = 62.07578
$ws2.Cells.Item(40, 3).Value = "forecast"
$ws2.Cells.Item(40, 5).Value = "OFF"
$ws2.Cells.Item(41, 1).Value = 46041.83333333334
$ws2.Cells.Item(41, 2).Value = 65
$ws2.Cells.Item(41, 3).Value = "forecast"
$ws2.Cells.Item(41, 5).Value = "ON"
$ws2.Cells.Item(42, 1).Value = 46041.85416666666
$ws2.Cells.Item(42, 2).Value = 65
$ws2.Cells.Item(42, 3).Value = "forecast"
$ws2.Cells.Item(42, 5).Value = "ON"
$ws2.Cells.Item(43, 1).Value = 46041.875
$ws2.Cells.Item(43, 2).Value = 65
$ws2.Cells.Item(43, 3).Value = "forecast"
$ws2.Cells.Item(43, 5).Value = "ON"
$ws2.Cells.Item(44, 1).Value = 46041.89583333334
$ws2.Cells.Item(44, 2).Value = 62.49071
$ws2.Cells.Item(44, 3).Value = "forecast"
$ws2.Cells.Item(44, 5).Value = "ON"
$ws2.Cells.Item(45, 1).Value = 46041.91666666666
$ws2.Cells.Item(45, 2).Value = 61.16151
$ws2.Cells.Item(45, 3).Value = "forecast"
$ws2.Cells.Item(45, 5).Value = "ON"
$ws2.Cells.Item(46, 1).Value = 46041.9375
$ws2.Cells.Item(46, 2).Value = 58.63887
$ws2.Cells.Item(46, 3).Value = "forecast"
$ws2.Cells.Item(46, 5).Value = "ON"
$ws2.Cells.Item(47, 1).Value = 46041.95833333334
$ws2.Cells.Item(47, 2).Value = 58.14989
$ws2.Cells.Item(47, 3).Value = "forecast"
$ws2.Cells.Item(47, 5).Value = "ON"
$ws2.Cells.Item(48, 1).Value = 46041.97916666666
$ws2.Cells.Item(48, 2).Value = 61.62715
$ws2.Cells.Item(48, 3).Value = "forecast"
$ws2.Cells.Item(48, 5).Value = "ON"

# Remove the now-superseded last row (old row 49)
$ws2.Rows.Item(49).Delete()
